$d = $word.ActiveDocument

# Replace the sentence, trimming everything from the en-dash onward
# and leaving a trailing space after "differences".
$find = [char]0x2013
$oldText = "sampling differences " + $find + " see manuscript and Gram-stain level results for more discussion."
$newText = "sampling differences "

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newText, 2)

# Remove the now-orphaned _GoBack bookmark, if present.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Save()
